# Fruta / hortaliza, semanal
#
# The weekly refresh re-shuffles the Fecha/Volumen/Precio columns among the
# existing rows (most other columns - Mercado, Region, Producto, Categoria,
# Variedad, Calidad, Unidad, Origen, Kg/unidad - stay put). Rows 10, 13, 14
# and 16 are not part of this week's refresh and keep their values.
#
# For each destination row below, D/M/N/O/P/S take on the values that used
# to live in the corresponding source row before the refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowValues {
    param($Row, $Fecha, $Volumen, $PrecioMin, $PrecioMax, $PrecioProm, $PrecioKg)

    $ws.Cells.Item($Row, 4).Value  = $Fecha       # D - Fecha
    $ws.Cells.Item($Row, 13).Value = $Volumen     # M - Volumen
    $ws.Cells.Item($Row, 14).Value = $PrecioMin   # N - Precio minimo
    $ws.Cells.Item($Row, 15).Value = $PrecioMax   # O - Precio maximo
    $ws.Cells.Item($Row, 16).Value = $PrecioProm  # P - Precio promedio ponderado
    $ws.Cells.Item($Row, 19).Value = $PrecioKg    # S - Precio $/Kg
}

# Row  2  <- old row 15
Set-RowValues 2  44490 400 9500  10000 9750  4875
# Row  3  <- old row 9
Set-RowValues 3  44881 440 6000  7000  6500  3250
# Row  4  <- old row 8
Set-RowValues 4  44818 200 11000 12000 11500 5750
# Row  5  <- old row 3
Set-RowValues 5  44874 300 7500  8000  7750  3875
# Row  6  <- old row 12
Set-RowValues 6  44889 460 3500  4000  3750  1875
# Row  7  <- old row 2
Set-RowValues 7  44454 160 12000 13000 12500 6250
# Row  8  <- old row 4
Set-RowValues 8  44875 400 7000  7500  7250  3625
# Row  9  <- old row 18
Set-RowValues 9  44517 400 5500  6000  5750  2875
# Row 11  <- old row 17
Set-RowValues 11 44895 240 3000  3500  3250  1625
# Row 12  <- old row 6
Set-RowValues 12 44819 240 11000 12000 11500 5750
# Row 15  <- old row 7
Set-RowValues 15 44489 160 9500  10000 9750  4875
# Row 17  <- old row 5
Set-RowValues 17 44497 500 9000  10000 9500  4750
# Row 18  <- old row 11
Set-RowValues 18 44455 200 12000 13000 12500 6250
